# South Carolina 2016 MCAS cleanup: rename headers to snake_case codes,
# title-case Spanish connector words (de/del/la/las/los/el/y) in the
# mx_state / mx_municipality columns, correct a couple of 1-ULP float
# artifacts in pct_matriculas, and drop the trailing footnote rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to short codes ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case connector words (de/del/la/los/las/y/el) in state/municipality names ---
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B25').Value = 'Amatenango De La Frontera'
$ws.Range('B26').Value = 'Amatenango Del Valle'
$ws.Range('B29').Value = 'Bejucal De Ocampo'
$ws.Range('B31').Value = 'Benemérito De Las Américas'
$ws.Range('B36').Value = 'Chiapa De Corzo'
$ws.Range('B41').Value = 'Comitán De Domínguez'
$ws.Range('B58').Value = 'Mazapa De Madero'
$ws.Range('B62').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B69').Value = 'Salto De Agua'
$ws.Range('B70').Value = 'San Cristóbal De Las Casas'
$ws.Range('B95').Value = 'Coyame Del Sotol'
$ws.Range('B99').Value = 'Guadalupe Y Calvo'
$ws.Range('B101').Value = 'Hidalgo Del Parral'
$ws.Range('B121').Value = 'San Juan De Sabinas'
$ws.Range('B130').Value = 'Villa De Álvarez'
$ws.Range('A132').Value = 'Ciudad De México'
$ws.Range('B155').Value = 'Nombre De Dios'
$ws.Range('B157').Value = 'Pánuco De Coronado'
$ws.Range('B160').Value = 'San Juan De Guadalupe'
$ws.Range('A165').Value = 'Estado De México'
$ws.Range('B165').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B168').Value = 'Almoloya De Alquisiras'
$ws.Range('B169').Value = 'Almoloya De Juárez'
$ws.Range('B173').Value = 'Atizapán De Zaragoza'
$ws.Range('B179').Value = 'Chapa De Mota'
$ws.Range('B182').Value = 'Coacalco De Berriozábal'
$ws.Range('B186').Value = 'Ecatepec De Morelos'
$ws.Range('B202').Value = 'Naucalpan De Juárez'
$ws.Range('B209').Value = 'San Felipe Del Progreso'
$ws.Range('B210').Value = 'San Martín De Las Pirámides'
$ws.Range('B211').Value = 'San Simón De Guerrero'
$ws.Range('B213').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B221').Value = 'Tenango Del Valle'
$ws.Range('B228').Value = 'Tlalnepantla De Baz'
$ws.Range('B233').Value = 'Valle De Bravo'
$ws.Range('B234').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B235').Value = 'Villa De Allende'
$ws.Range('B236').Value = 'Villa Del Carbón'
$ws.Range('B246').Value = 'San Miguel De Allende'
$ws.Range('B247').Value = 'Apaseo El Alto'
$ws.Range('B248').Value = 'Apaseo El Grande'
$ws.Range('B256').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B260').Value = 'Jaral Del Progreso'
$ws.Range('B267').Value = 'Purísima Del Rincón'
$ws.Range('B271').Value = 'San Diego De La Unión'
$ws.Range('B273').Value = 'San Francisco Del Rincón'
$ws.Range('B275').Value = 'San Luis De La Paz'
$ws.Range('B277').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B278').Value = 'Silao De La Victoria'
$ws.Range('B282').Value = 'Valle De Santiago'
$ws.Range('B288').Value = 'Acapulco De Juárez'
$ws.Range('B290').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B291').Value = 'Alcozauca De Guerrero'
$ws.Range('B294').Value = 'Atenango Del Río'
$ws.Range('B296').Value = 'Atoyac De Álvarez'
$ws.Range('B297').Value = 'Ayutla De Los Libres'
$ws.Range('B300').Value = 'Buenavista De Cuéllar'
$ws.Range('B301').Value = 'Chilapa De Álvarez'
$ws.Range('B302').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B303').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B306').Value = 'Coyuca De Benítez'
$ws.Range('B307').Value = 'Coyuca De Catalán'
$ws.Range('B311').Value = 'Cutzamala De Pinzón'
$ws.Range('B316').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B317').Value = 'Iguala De La Independencia'
$ws.Range('B319').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B320').Value = 'Zihuatanejo De Azueta'
$ws.Range('B322').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B325').Value = 'Mártir De Cuilapan'
$ws.Range('B337').Value = 'Taxco De Alarcón'
$ws.Range('B339').Value = 'Técpan De Galeana'
$ws.Range('B341').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B343').Value = 'Tixtla De Guerrero'
$ws.Range('B347').Value = 'Tlapa De Comonfort'
$ws.Range('B359').Value = 'Agua Blanca De Iturbide'
$ws.Range('B365').Value = 'Atotonilco El Grande'
$ws.Range('B371').Value = 'Cuautepec De Hinojosa'
$ws.Range('B374').Value = 'Huasca De Ocampo'
$ws.Range('B377').Value = 'Huejutla De Reyes'
$ws.Range('B380').Value = 'Jacala De Ledezma'
$ws.Range('B385').Value = 'Mineral Del Chico'
$ws.Range('B386').Value = 'Mineral Del Monte'
$ws.Range('B387').Value = 'Molango De Escamilla'
$ws.Range('B389').Value = 'Nopala De Villagrán'
$ws.Range('B390').Value = 'Omitlán De Juárez'
$ws.Range('B391').Value = 'Pachuca De Soto'
$ws.Range('B394').Value = 'Progreso De Obregón'
$ws.Range('B399').Value = 'Santiago De Anaya'
$ws.Range('B400').Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range('B404').Value = 'Tenango De Doria'
$ws.Range('B405').Value = 'Tepehuacán De Guerrero'
$ws.Range('B406').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B407').Value = 'Tezontepec De Aldama'
$ws.Range('B413').Value = 'Tula De Allende'
$ws.Range('B414').Value = 'Tulancingo De Bravo'
$ws.Range('B417').Value = 'Zacualtipán De Ángeles'
$ws.Range('B422').Value = 'Ahualulco De Mercado'
$ws.Range('B434').Value = 'Encarnación De Díaz'
$ws.Range('B438').Value = 'Ixtlahuacán Del Río'
$ws.Range('B441').Value = 'Jilotlán De Los Dolores'
$ws.Range('B444').Value = 'La Manzanilla De La Paz'
$ws.Range('B445').Value = 'Lagos De Moreno'
$ws.Range('B448').Value = 'Ojuelos De Jalisco'
$ws.Range('B453').Value = 'San Juan De Los Lagos'
$ws.Range('B454').Value = 'San Miguel El Alto'
$ws.Range('B456').Value = 'Talpa De Allende'
$ws.Range('B457').Value = 'Tamazula De Gordiano'
$ws.Range('B461').Value = 'Tepatitlán De Morelos'
$ws.Range('B463').Value = 'Tizapán El Alto'
$ws.Range('B464').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B468').Value = 'Unión De San Antonio'
$ws.Range('B469').Value = 'Unión De Tula'
$ws.Range('B471').Value = 'Yahualica De González Gallo'
$ws.Range('B473').Value = 'Zapotlán El Grande'
$ws.Range('B490').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B492').Value = 'Cojumatlán De Régules'
$ws.Range('B536').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B556').Value = 'Coatlán Del Río'
$ws.Range('B563').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B567').Value = 'Puente De Ixtla'
$ws.Range('B571').Value = 'Tetela Del Volcán'
$ws.Range('B581').Value = 'Bahía De Banderas'
$ws.Range('B584').Value = 'Ixtlán Del Río'
$ws.Range('B590').Value = 'Santa María Del Oro'
$ws.Range('B609').Value = 'San Nicolás De Los Garza'
$ws.Range('B614').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B620').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B622').Value = 'Coicoyán De Las Flores'
$ws.Range('B623').Value = 'Constancia Del Rosario'
$ws.Range('B625').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B626').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B627').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B628').Value = 'Ixtlán De Juárez'
$ws.Range('B629').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B635').Value = 'Mártires De Tacubaya'
$ws.Range('B637').Value = 'Mazatlán Villa De Flores'
$ws.Range('B638').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B639').Value = 'Nejapa De Madero'
$ws.Range('B640').Value = 'Oaxaca De Juárez'
$ws.Range('B641').Value = 'Ocotlán De Morelos'
$ws.Range('B642').Value = 'Pinotepa De Don Luis'
$ws.Range('B643').Value = 'Putla Villa De Guerrero'
$ws.Range('B644').Value = 'Reforma De Pineda'
$ws.Range('B657').Value = 'San Dionisio Del Mar'
$ws.Range('B660').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B669').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B698').Value = 'San Miguel Del Puerto'
$ws.Range('B702').Value = 'San Pablo Villa De Mitla'
$ws.Range('B719').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B722').Value = 'Santa Inés Del Monte'
$ws.Range('B731').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B755').Value = 'Santo Domingo De Morelos'
$ws.Range('B763').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B764').Value = 'Tataltepec De Valdés'
$ws.Range('B765').Value = 'Teotitlán De Flores Magón'
$ws.Range('B766').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B767').Value = 'Tlacolula De Matamoros'
$ws.Range('B769').Value = 'Villa De Tututepec'
$ws.Range('B770').Value = 'Villa De Zaachila'
$ws.Range('B772').Value = 'Villa Sola De Vega'
$ws.Range('B773').Value = 'Zapotitlán Del Río'
$ws.Range('B774').Value = 'Zimatlán De Álvarez'
$ws.Range('B786').Value = 'Ayotoxco De Guerrero'
$ws.Range('B788').Value = 'Chalchicomula De Sesma'
$ws.Range('B796').Value = 'Chila De La Sal'
$ws.Range('B804').Value = 'Cuetzalan Del Progreso'
$ws.Range('B816').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B818').Value = 'Izúcar De Matamoros'
$ws.Range('B824').Value = 'Mazapiltepec De Juárez'
$ws.Range('B829').Value = 'Palmar De Bravo'
$ws.Range('B844').Value = 'San Salvador El Seco'
$ws.Range('B853').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B856').Value = 'Tepexi De Rodríguez'
$ws.Range('B858').Value = 'Tetela De Ocampo'
$ws.Range('B862').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B870').Value = 'Xayacatlán De Bravo'
$ws.Range('B878').Value = 'Zapotitlán De Méndez'
$ws.Range('B883').Value = 'Amealco De Bonfil'
$ws.Range('B885').Value = 'Cadereyta De Montes'
$ws.Range('B888').Value = 'Jalpan De Serra'
$ws.Range('B889').Value = 'Landa De Matamoros'
$ws.Range('B892').Value = 'Pinal De Amoles'
$ws.Range('B895').Value = 'San Juan Del Río'
$ws.Range('B905').Value = 'Armadillo De Los Infante'
$ws.Range('B906').Value = 'Axtla De Terrazas'
$ws.Range('B910').Value = 'Ciudad Del Maíz'
$ws.Range('B920').Value = 'Mexquitic De Carmona'
$ws.Range('B924').Value = 'San Ciro De Acosta'
$ws.Range('B929').Value = 'Santa María Del Río'
$ws.Range('B938').Value = 'Villa De Arista'
$ws.Range('B939').Value = 'Villa De Guadalupe'
$ws.Range('B940').Value = 'Villa De Ramos'
$ws.Range('B941').Value = 'Villa De Reyes'
$ws.Range('B963').Value = 'Nacozari De García'
$ws.Range('B975').Value = 'Jalpa De Méndez'
$ws.Range('B995').Value = 'Soto La Marina'
$ws.Range('B1007').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B1010').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B1011').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B1014').Value = 'San Pablo Del Monte'
$ws.Range('B1027').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B1031').Value = 'Amatlán De Los Reyes'
$ws.Range('B1041').Value = 'Boca Del Río'
$ws.Range('B1043').Value = 'Camarón De Tejeda'
$ws.Range('B1061').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1062').Value = 'Cosautlán De Carvajal'
$ws.Range('B1077').Value = 'Hueyapan De Ocampo'
$ws.Range('B1078').Value = 'Ignacio De La Llave'
$ws.Range('B1081').Value = 'Ixhuacán De Los Reyes'
$ws.Range('B1082').Value = 'Ixhuatlán De Madero'
$ws.Range('B1083').Value = 'Ixhuatlán Del Café'
$ws.Range('B1084').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1091').Value = 'Juchique De Ferrer'
$ws.Range('B1095').Value = 'Lerdo De Tejada'
$ws.Range('B1099').Value = 'Martínez De La Torre'
$ws.Range('B1102').Value = 'Medellín De Bravo'
$ws.Range('B1105').Value = 'Mixtla De Altamirano'
$ws.Range('B1114').Value = 'Ozuluama De Mascareñas'
$ws.Range('B1117').Value = 'Paso De Ovejas'
$ws.Range('B1118').Value = 'Paso Del Macho'
$ws.Range('B1122').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1130').Value = 'Sayula De Alemán'
$ws.Range('B1134').Value = 'Soledad De Doblado'
$ws.Range('B1137').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1163').Value = 'Vega De Alatorre'
$ws.Range('B1179').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1188').Value = 'Jiménez Del Teul'
$ws.Range('B1190').Value = 'Mezquital Del Oro'
$ws.Range('B1192').Value = 'Nochistlán De Mejía'
$ws.Range('B1193').Value = 'Noria De Ángeles'
$ws.Range('B1198').Value = 'Teúl De González Ortega'
$ws.Range('B1201').Value = 'Villa De Cos'

# --- Floating point recalculation artifacts on percentage column (1 ULP) ---
$ws.Range('D25').Value = 0.0009068923821039904
$ws.Range('D206').Value = 0.0009068923821039904
$ws.Range('D329').Value = 0.0009068923821039904
$ws.Range('D354').Value = 0.0009068923821039904
$ws.Range('D634').Value = 0.0009068923821039904
$ws.Range('D658').Value = 0.0009068923821039904
$ws.Range('D664').Value = 0.0009068923821039904
$ws.Range('D678').Value = 0.0009068923821039904
$ws.Range('D717').Value = 0.0009068923821039904
$ws.Range('D727').Value = 0.0009068923821039904
$ws.Range('D761').Value = 0.0009068923821039904
$ws.Range('D783').Value = 0.0009068923821039904
$ws.Range('D815').Value = 0.0009068923821039904
$ws.Range('D818').Value = 0.0009068923821039904
$ws.Range('D828').Value = 0.0009068923821039904
$ws.Range('D888').Value = 0.0009068923821039904
$ws.Range('D900').Value = 0.0009068923821039904
$ws.Range('D924').Value = 0.0009068923821039904
$ws.Range('D929').Value = 0.0009068923821039904
$ws.Range('D1009').Value = 0.0009068923821039904
$ws.Range('D1066').Value = 0.0009068923821039904
$ws.Range('D1067').Value = 0.0009068923821039904
$ws.Range('D1095').Value = 0.0009068923821039904
$ws.Range('D1102').Value = 0.0009068923821039904
$ws.Range('D1146').Value = 0.0009068923821039904
$ws.Range('D1168').Value = 0.0009068923821039904

# --- Remove trailing footnote/metadata rows; row 1206 (grand Total) is now last ---
$ws.Range("A1208:A1212").EntireRow.Delete()
